$d = $word.ActiveDocument

# --- Change 1: paragraph "O Pao quando nao e partido..." -> accent the
#     first "e" ("e" -> "e" with acute accent) and in doing so the run
#     that held the whole sentence gets split into three runs:
#       "O Pao quando nao "  |  "e" (accented)  |  " partido ele representa..."
$searchRange = $d.Content.Duplicate
$needle = "n" + [char]0x00E3 + "o e partido"
$found = $searchRange.Find.Execute($needle)

# Position of the bare "e" that needs the acute accent: right after the
# leading "nao " (4 characters) inside the matched snippet.
$eStart = $searchRange.Start + 4
$eRange = $d.Range($eStart, $eStart + 1)

# Toggling a character-level property on just this one-letter range forces
# the engine to split the run into three (matching how Word itself splits
# runs at an edit boundary); then we flip the property back off and set
# the corrected, accented text so the run keeps the paragraph's normal
# (non-bold) formatting while remaining its own run.
$eRange.Bold = 1
$eRange.Text = [char]0x00E9
$eRange2 = $d.Range($eStart, $eStart + 1)
$eRange2.Bold = 0

# --- Change 2: fix the typo "Dues" -> "Deus" inside the first comment.
$comment = $d.Comments(1)
$commentRange = $comment.Range.Duplicate
$commentText = $commentRange.Text
$fixedText = $commentText.Replace("ser Dues", "ser Deus")
if ($fixedText -eq $commentText) {
    $fixedText = $commentText.Replace("Dues", "Deus")
}
$commentRange.Text = $fixedText
